$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-applied kmeans clustering results (5 clusters) -- corrected row data for B2:N21
$data = @(
    @('Europe', 1.36, 19840.3, 0.89, 0, 77.67, 127.08, 28.19, 100, 0, 3, 4, 1),
    @('Europe', 1.31, 35911.7, 0.96, 7.689999999999998, 89, 110.12, 26.75, 100, -9, 3, 4, 1),
    @('Europe', 0.09, 11341.3, 0.99, 0, 58.64, 93.35, 24.88, 100, 1, 3, 4, 1),
    @('Europe', 0.89, 26980.3, 0.83, 0, 64.35, 99.72, 38.86, 100, 0, 3, 4, 1),
    @('Europe', 0.12, 13102.1, 0.86, -12.5, 55.29, 99.4, 50.91, 100, -2, 3, 5, 1),
    @('Asie', 0.68, 59842.2, 0.98, -223.93, 69.41, 192.94, 23.26, 100, 24, 1, 4, 1),
    @('Océanie', 2.34, 126144, 0.85, -21.05, 67.75, 137.1, 17.5, 100, -16, 1, 4, 1),
    @('Asie', 2, 11366.3, 0.64, -11.11, 77.72, 129.03, 1.48, 100, -2, 3, 4, 1),
    @('Europe', 0.17, 53071.5, 0.82, 51.58, 60.54, 127.02, 12.63, 100, 57, 1, 4, 1),
    @('Amériques', 0.61, 54173, 1.2, 36.42, 59.85, 130.44, 11.58, 100, -16, 1, 4, 1),
    @('Amériques', 0.68, 50442.3, 0.73, -122.37, 56.93, 116.33, 7.92, 100, 13, 1, 4, 1),
    @('Amériques', 0.48, 55356.7, 0.98, 20.36, 68.19, 155.24, 16.84, 100, 0, 1, 4, 1),
    @('Europe', 0.99, 77749.2, 1.01, 22.66, 64.27, 137.1, 17.06, 100, 24, 1, 4, 1),
    @('Europe', 2.46, 114986, 1.41, 0, 63.9, 136.29, 10.64, 100, 1, 1, 4, 1),
    @('Amériques', 0.74, 43559.3, 1.16, 35.70999999999999, 55.62, 130.44, 19.6, 100, 2, 3, 4, 1),
    @('Europe', 0.31, 55088.6, 0.98, -63.44, 65.12, 137.1, 12.47, 100, 164, 1, 4, 1),
    @('Europe', 1.22, 69103.6, 1.31, 61.65, 63.59, 122.38, 10.69, 100, -10, 1, 4, 1),
    @('Océanie', 1.39, 31100, 0.45, 18.18000000000001, 60.51, 101.41, 27.93, 100, 2, 3, 4, 1),
    @('Océanie', 0.51, 17000, -0.07, 0, 63.89, 125.6, 26, 100, -1, 3, 4, 1),
    @('Asie', 0.73, 6486.1, 1.04, -13.33, 61.03, 105.44, 41.84, 100, -4, 3, 5, 1)
)

$colLetters = @("B","C","D","E","F","G","H","I","J","K","L","M","N")

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = 2 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 2).Value = $rowValues[$j]
    }
}

Write-Host "Updated rows 2-21, columns B-N"